$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.807.89"
$ws.Range("E2").Value = "  +3.15%  "
$ws.Range("D3").Value = "1.866.65"
$ws.Range("E3").Value = "  +2.82%  "
$ws.Range("E4").Value = "  +3.31%  "
$ws.Range("D5").Value = "324.90"
$ws.Range("E5").Value = "  +4.03%  "
$ws.Range("D6").Value = "1.036"
$ws.Range("E6").Value = "  +3.03%  "
$ws.Range("D7").Value = "0.4426"
$ws.Range("E7").Value = "  +2.94%  "
$ws.Range("D8").Value = "0.3803"
$ws.Range("E8").Value = "  +3.63%  "
$ws.Range("D9").Value = "0.07486"
$ws.Range("E9").Value = "  +3.13%  "
$ws.Range("D10").Value = "0.8861"
$ws.Range("E10").Value = "  +2.55%  "
$ws.Range("D11").Value = "21.78"
$ws.Range("E11").Value = "  +2.31%  "
$ws.Range("D12").Value = "1.888.36"
$ws.Range("E12").Value = "  -12.39%  "
$ws.Range("D13").Value = "5.566"
$ws.Range("E13").Value = "  +2.90%  "
$ws.Range("D14").Value = "6.765"
$ws.Range("E14").Value = "  +2.44%  "
$ws.Range("D15").Value = "0.07237"
$ws.Range("E15").Value = "  +4.26%  "
$ws.Range("D16").Value = "83.94"
$ws.Range("E16").Value = "  +3.64%  "
$ws.Range("D17").Value = "1.041"
$ws.Range("E17").Value = "  +2.89%  "
$ws.Range("D18").Value = "0.000009168"
$ws.Range("E18").Value = "  +3.42%  "
$ws.Range("D19").Value = "1.037"
$ws.Range("D20").Value = "15.58"
$ws.Range("E20").Value = "  +2.12%  "
$ws.Range("D21").Value = "27.815.87"
$ws.Range("E21").Value = "  +2.98%  "
$ws.Range("D22").Value = "5.331"
$ws.Range("E22").Value = "  +2.81%  "
$ws.Range("D23").Value = "11.39"
$ws.Range("E23").Value = "  +3.56%  "
$ws.Range("D24").Value = "1.999"
$ws.Range("E24").Value = "  +6.09%  "
$ws.Range("D25").Value = "158.83"
$ws.Range("E25").Value = "  +3.17%  "
$ws.Range("D26").Value = "18.91"
$ws.Range("E26").Value = "  +3.34%  "
$ws.Range("D27").Value = "5.340"
$ws.Range("E27").Value = "  +2.30%  "
$ws.Range("D28").Value = "1.988"
$ws.Range("E28").Value = "  +4.71%  "
$ws.Range("D29").Value = "117.87"
$ws.Range("E29").Value = "  +2.83%  "
$ws.Range("D30").Value = "0.09083"
$ws.Range("E30").Value = "  +1.59%  "
$ws.Range("D31").Value = "3.131"
$ws.Range("E31").Value = "  +11.42%  "
$ws.Range("D32").Value = "0.7805"
$ws.Range("E32").Value = "  +4.62%  "
$ws.Range("D33").Value = "1.217"
$ws.Range("E33").Value = "  +2.70%  "
$ws.Range("D34").Value = "4.585"
$ws.Range("E34").Value = "  +3.85%  "
$ws.Range("E35").Value = "  +3.22%  "
$ws.Range("D36").Value = "1.156"
$ws.Range("E36").Value = "  +2.27%  "
$ws.Range("D37").Value = "0.01997"
$ws.Range("E37").Value = "  +3.81%  "
$ws.Range("D38").Value = "0.05360"
$ws.Range("E38").Value = "  +2.82%  "
$ws.Range("D39").Value = "2.872"
$ws.Range("E39").Value = "  +4.96%  "
$ws.Range("D40").Value = "0.5215"
$ws.Range("E40").Value = "  +2.31%  "
$ws.Range("D41").Value = "0.1698"
$ws.Range("E41").Value = "  +2.61%  "
$ws.Range("D42").Value = "6.944"
$ws.Range("E42").Value = "  +7.17%  "
$ws.Range("D43").Value = "8.699"
$ws.Range("E43").Value = "  +4.41%  "
$ws.Range("D44").Value = "10.78"
$ws.Range("E44").Value = "  +3.45%  "
$ws.Range("D45").Value = "109.74"
$ws.Range("E45").Value = "  +2.85%  "
$ws.Range("D46").Value = "1.730"
$ws.Range("E46").Value = "  +5.32%  "
$ws.Range("D47").Value = "0.4724"
$ws.Range("E47").Value = "  +3.03%  "
$ws.Range("D48").Value = "0.06472"
$ws.Range("E48").Value = "  +4.12%  "
$ws.Range("D49").Value = "1.911"
$ws.Range("E49").Value = "  +4.17%  "
$ws.Range("D50").Value = "39.95"
$ws.Range("E50").Value = "  +4.07%  "
$ws.Range("D51").Value = "64.70"
$ws.Range("E51").Value = "  +2.67%  "
